# "Refactoring.. - updated saek"
#
# Core edit: cell AG2 on the "epas" sheet held the number 1; it is
# changed to the (new) text value "er". Everything else in the target
# diff (shared-string table re-numbering, AH2/AI2/AJ2 shared-string
# index bumps) is a mechanical side effect of inserting that new shared
# string ahead of the existing ones and is handled automatically by the
# engine when the cell value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AG2").Value = "er"

# Reflect the saved view state from the diff: the window was scrolled so
# column X is the left-most visible column, and the selection moved from
# C2 to the cell that was just edited (AG2).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 24
$ws.Range("AG2").Select()

# Tiny column-C width touch-up (38.32 -> 38.33 characters) recorded in
# the diff alongside the edit.
$ws.Columns.Item(3).ColumnWidth = 37.5
